# datapenjualan_bushub.xlsx
# perbaiki logika bayar di routeTransaction.py bagian insert detailTransaksi
#
# Update the values of rows 5-8 to reflect the corrected payment logic,
# and remove the old row 9 (A0010) which no longer exists after the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain TEXT (the source data stores
# numbers/dates as strings, e.g. "2280000.0" or "2024-11-25") even though it
# looks numeric/date-like. Force the cell to Text format first so Excel does
# not auto-convert it to a number/date, then drop the now-unneeded
# formatting so the cell keeps the workbook's default (unstyled) look.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Row 5 (A0003) ---
$ws.Range("B5").Value = "2024-11-25 18:03:54.803392+07:00"
$ws.Range("C5").Value = "user2@gmail.com"
$ws.Range("D5").Value = "KSR_01"
$ws.Range("E5").Value = "COMPLETED"
Set-TextValue $ws.Range("F5") "2280000.0"
$ws.Range("G5").Value = "transfer"
$ws.Range("K5").Value = ""

# --- Row 6 (A0004) ---
$ws.Range("B6").Value = "2024-11-25 18:05:21.359478+07:00"
$ws.Range("C6").Value = "user2@gmail.com"
$ws.Range("D6").Value = "KSR_01"
Set-TextValue $ws.Range("F6") "500000.0"
Set-TextValue $ws.Range("H6") "2024-11-25"
Set-TextValue $ws.Range("I6") "2024-11-26"
$ws.Range("K6").Value = "Paket Wisata Singkawang"

# --- Row 7 (was A0006, becomes A0003) ---
$ws.Range("A7").Value = "A0003"
$ws.Range("B7").Value = "2024-11-25 18:03:54.803392+07:00"
$ws.Range("C7").Value = "user2@gmail.com"
$ws.Range("D7").Value = "KSR_01"
$ws.Range("E7").Value = "COMPLETED"
Set-TextValue $ws.Range("F7") "2280000.0"
Set-TextValue $ws.Range("H7") "2024-11-25"
Set-TextValue $ws.Range("I7") "2024-11-28"

# --- Row 8 (was A0007, becomes A0004) ---
$ws.Range("A8").Value = "A0004"
$ws.Range("B8").Value = "2024-11-25 18:05:21.359478+07:00"
$ws.Range("C8").Value = "user2@gmail.com"
$ws.Range("D8").Value = "KSR_01"

# --- Row 9 (A0010) no longer exists; remove it entirely ---
$ws.Rows.Item(9).Delete()
